{"js": "// Convert the Word field { m:null.fit(null, null) } (fldChar begin/instrText/fldChar end)\n// into literal M2Doc-style template text \"{m:null.fit(null, null)}\" made of plain <w:t> runs,\n// keeping the orange color run-formatting on the \"null.fit(null, null)\" portion.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the paragraph that contains the field (fldChar begin/end, no plain text).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const fields = para.fields;\n  fields.load(\"items\");\n  await context.sync();\n  if (fields.items.length > 0) {\n    const f = fields.items[0];\n    f.load(\"code\");\n    await context.sync();\n    if (f.code.indexOf(\"null.fit(\") !== -1) {\n      target = para;\n      break;\n    }\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the target field paragraph.\");\n}\n\n// Read back this paragraph's own OOXML so we can keep its exact <w:p> attributes\n// (w14:paraId, rsids, ...) untouched - the diff only rewrites the runs inside it.\nconst ooxmlProp = target.getOoxml();\nawait context.sync();\nconst fullOoxml = ooxmlProp.value;\n\nconst pOpenMatch = fullOoxml.match(/<w:p(?:\\s[^>]*)?>/);\nconst pOpenTag = pOpenMatch ? pOpenMatch[0] : \"<w:p>\";\n\nconst colorRpr = '<w:rPr><w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/></w:rPr>';\n\nconst newRunsXml =\n  '<w:r><w:t>{</w:t></w:r>' +\n  '<w:r><w:t>m</w:t></w:r>' +\n  '<w:r><w:t>:</w:t></w:r>' +\n  '<w:r>' + colorRpr + '<w:t>null</w:t></w:r>' +\n  '<w:r>' + colorRpr + '<w:t>.fit(</w:t></w:r>' +\n  '<w:r>' + colorRpr + '<w:t>null</w:t></w:r>' +\n  '<w:r>' + colorRpr + '<w:t xml:space=\"preserve\">, </w:t></w:r>' +\n  '<w:r>' + colorRpr + '<w:t>null</w:t></w:r>' +\n  '<w:r>' + colorRpr + '<w:t>)</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>';\n\nconst replacementOoxml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n  '<w:body>' + pOpenTag + newRunsXml + '</w:p></w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nconst fullRange = target.getRange();\nfullRange.insertOoxml(replacementOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Convert the Word field { m:null.fit(null, null) } (fldChar begin/instrText/fldChar end)\n# into literal M2Doc-style template text \"{m:null.fit(null, null)}\" made of plain text runs,\n# keeping the orange color run-formatting on the \"null.fit(null, null)\" portion.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that holds the field (fldChar begin/instrText/fldChar end).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Fields.Count -gt 0) {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not find the paragraph containing the field.\"\n}\n\n# Remove the field itself (fldChar begin/instrText.../fldChar end runs); this leaves\n# an empty paragraph (just the paragraph mark) ready to receive the literal text runs.\n$fld = $target.Range.Fields.Item(1)\n$fld.Delete()\n\n# Collapsed range at the (now empty) start of that paragraph's content - inserting a\n# <w:p> fragment here merges the new runs into the existing paragraph and keeps the\n# paragraph's own attributes (w14:paraId, rsids, ...) untouched.\n$insertPos = $target.Range.Start\n$rng = $d.Range($insertPos, $insertPos)\n\n$colorRPr = \"<w:rPr><w:color w:val='E36C0A' w:themeColor='accent6' w:themeShade='BF'/></w:rPr>\"\n\n$runsXml = \"<w:r><w:t>{</w:t></w:r>\" + `\n           \"<w:r><w:t>m</w:t></w:r>\" + `\n           \"<w:r><w:t>:</w:t></w:r>\" + `\n           \"<w:r>$colorRPr<w:t>null</w:t></w:r>\" + `\n           \"<w:r>$colorRPr<w:t>.fit(</w:t></w:r>\" + `\n           \"<w:r>$colorRPr<w:t>null</w:t></w:r>\" + `\n           \"<w:r>$colorRPr<w:t xml:space='preserve'>, </w:t></w:r>\" + `\n           \"<w:r>$colorRPr<w:t>null</w:t></w:r>\" + `\n           \"<w:r>$colorRPr<w:t>)</w:t></w:r>\" + `\n           \"<w:r><w:t xml:space='preserve'>}</w:t></w:r>\"\n\n$xmlFrag = \"<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>$runsXml</w:p>\"\n\n$rng.InsertXML($xmlFrag)\n"}
